$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-05-12 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-13 Saturday", 2) | Out-Null
$d.Content.Find.Execute("85×64=", $true, $false, $false, $false, $false, $true, 1, $false, "19×57=", 2) | Out-Null
$d.Content.Find.Execute("59×60=", $true, $false, $false, $false, $false, $true, 1, $false, "59×27=", 2) | Out-Null
$d.Content.Find.Execute("81×24=", $true, $false, $false, $false, $false, $true, 1, $false, "90×30=", 2) | Out-Null
$d.Content.Find.Execute("21×37=", $true, $false, $false, $false, $false, $true, 1, $false, "70×59=", 2) | Out-Null
$d.Content.Find.Execute("63×74=", $true, $false, $false, $false, $false, $true, 1, $false, "93×24=", 2) | Out-Null
$d.Content.Find.Execute("44×62=", $true, $false, $false, $false, $false, $true, 1, $false, "12×91=", 2) | Out-Null
$d.Content.Find.Execute("13×42=", $true, $false, $false, $false, $false, $true, 1, $false, "32×60=", 2) | Out-Null
$d.Content.Find.Execute("46×17=", $true, $false, $false, $false, $false, $true, 1, $false, "84×46=", 2) | Out-Null
$d.Content.Find.Execute("28×41=", $true, $false, $false, $false, $false, $true, 1, $false, "87×30=", 2) | Out-Null
$d.Content.Find.Execute("57×86=", $true, $false, $false, $false, $false, $true, 1, $false, "98×28=", 2) | Out-Null
$d.Content.Find.Execute("24×61=", $true, $false, $false, $false, $false, $true, 1, $false, "56×92=", 2) | Out-Null
$d.Content.Find.Execute("54×35=", $true, $false, $false, $false, $false, $true, 1, $false, "73×44=", 2) | Out-Null
$d.Content.Find.Execute("68×73=", $true, $false, $false, $false, $false, $true, 1, $false, "92×64=", 2) | Out-Null
$d.Content.Find.Execute("97×88=", $true, $false, $false, $false, $false, $true, 1, $false, "36×97=", 2) | Out-Null
$d.Content.Find.Execute("67×79=", $true, $false, $false, $false, $false, $true, 1, $false, "82×63=", 2) | Out-Null
$d.Content.Find.Execute("77×66=", $true, $false, $false, $false, $false, $true, 1, $false, "13×84=", 2) | Out-Null
$d.Content.Find.Execute("87×54=", $true, $false, $false, $false, $false, $true, 1, $false, "85×51=", 2) | Out-Null
$d.Content.Find.Execute("79×30=", $true, $false, $false, $false, $false, $true, 1, $false, "97×52=", 2) | Out-Null
$d.Content.Find.Execute("85×91=", $true, $false, $false, $false, $false, $true, 1, $false, "35×82=", 2) | Out-Null
$d.Content.Find.Execute("78×63=", $true, $false, $false, $false, $false, $true, 1, $false, "90×100=", 2) | Out-Null
$d.Content.Find.Execute("75×82=", $true, $false, $false, $false, $false, $true, 1, $false, "66×18=", 2) | Out-Null
$d.Content.Find.Execute("47×26=", $true, $false, $false, $false, $false, $true, 1, $false, "14×28=", 2) | Out-Null
$d.Content.Find.Execute("26×21=", $true, $false, $false, $false, $false, $true, 1, $false, "69×69=", 2) | Out-Null
$d.Content.Find.Execute("16×57=", $true, $false, $false, $false, $false, $true, 1, $false, "32×69=", 2) | Out-Null
$d.Content.Find.Execute("49×16=", $true, $false, $false, $false, $false, $true, 1, $false, "53×48=", 2) | Out-Null
$d.Content.Find.Execute("97×57=", $true, $false, $false, $false, $false, $true, 1, $false, "50×62=", 2) | Out-Null
$d.Content.Find.Execute("85×60=", $true, $false, $false, $false, $false, $true, 1, $false, "29×57=", 2) | Out-Null
$d.Content.Find.Execute("43×79=", $true, $false, $false, $false, $false, $true, 1, $false, "88×18=", 2) | Out-Null
$d.Content.Find.Execute("28×34=", $true, $false, $false, $false, $false, $true, 1, $false, "67×33=", 2) | Out-Null
$d.Content.Find.Execute("35×46=", $true, $false, $false, $false, $false, $true, 1, $false, "83×21=", 2) | Out-Null
$d.Content.Find.Execute("56×12=", $true, $false, $false, $false, $false, $true, 1, $false, "54×11=", 2) | Out-Null
$d.Content.Find.Execute("73×89=", $true, $false, $false, $false, $false, $true, 1, $false, "81×51=", 2) | Out-Null
$d.Content.Find.Execute("78×52=", $true, $false, $false, $false, $false, $true, 1, $false, "60×21=", 2) | Out-Null
$d.Content.Find.Execute("65×15=", $true, $false, $false, $false, $false, $true, 1, $false, "96×93=", 2) | Out-Null
$d.Content.Find.Execute("55×31=", $true, $false, $false, $false, $false, $true, 1, $false, "67×83=", 2) | Out-Null
$d.Content.Find.Execute("18×63=", $true, $false, $false, $false, $false, $true, 1, $false, "30×12=", 2) | Out-Null
$d.Content.Find.Execute("25×15=", $true, $false, $false, $false, $false, $true, 1, $false, "72×12=", 2) | Out-Null
$d.Content.Find.Execute("96×97=", $true, $false, $false, $false, $false, $true, 1, $false, "57×30=", 2) | Out-Null
$d.Content.Find.Execute("21×33=", $true, $false, $false, $false, $false, $true, 1, $false, "94×36=", 2) | Out-Null
$d.Content.Find.Execute("66×28=", $true, $false, $false, $false, $false, $true, 1, $false, "32×20=", 2) | Out-Null
$d.Content.Find.Execute("89×10=", $true, $false, $false, $false, $false, $true, 1, $false, "21×72=", 2) | Out-Null
$d.Content.Find.Execute("81×73=", $true, $false, $false, $false, $false, $true, 1, $false, "46×96=", 2) | Out-Null
$d.Content.Find.Execute("34×87=", $true, $false, $false, $false, $false, $true, 1, $false, "21×30=", 2) | Out-Null
$d.Content.Find.Execute("73×98=", $true, $false, $false, $false, $false, $true, 1, $false, "66×38=", 2) | Out-Null
$d.Content.Find.Execute("81×21=", $true, $false, $false, $false, $false, $true, 1, $false, "44×52=", 2) | Out-Null
$d.Content.Find.Execute("20×60=", $true, $false, $false, $false, $false, $true, 1, $false, "25×49=", 2) | Out-Null
$d.Content.Find.Execute("66×22=", $true, $false, $false, $false, $false, $true, 1, $false, "62×31=", 2) | Out-Null
$d.Content.Find.Execute("15×75=", $true, $false, $false, $false, $false, $true, 1, $false, "68×29=", 2) | Out-Null
$d.Content.Find.Execute("73×58=", $true, $false, $false, $false, $false, $true, 1, $false, "22×28=", 2) | Out-Null
$d.Content.Find.Execute("40×99=", $true, $false, $false, $false, $false, $true, 1, $false, "11×65=", 2) | Out-Null
$d.Content.Find.Execute("70×55=", $true, $false, $false, $false, $false, $true, 1, $false, "30×99=", 2) | Out-Null
$d.Content.Find.Execute("52×81=", $true, $false, $false, $false, $false, $true, 1, $false, "97×92=", 2) | Out-Null
$d.Content.Find.Execute("76×88=", $true, $false, $false, $false, $false, $true, 1, $false, "11×80=", 2) | Out-Null
$d.Content.Find.Execute("27×76=", $true, $false, $false, $false, $false, $true, 1, $false, "72×93=", 2) | Out-Null
$d.Content.Find.Execute("74×42=", $true, $false, $false, $false, $false, $true, 1, $false, "92×21=", 2) | Out-Null
$d.Content.Find.Execute("35×44=", $true, $false, $false, $false, $false, $true, 1, $false, "24×92=", 2) | Out-Null
$d.Content.Find.Execute("42×88=", $true, $false, $false, $false, $false, $true, 1, $false, "78×20=", 2) | Out-Null
$d.Content.Find.Execute("81×70=", $true, $false, $false, $false, $false, $true, 1, $false, "17×55=", 2) | Out-Null
$d.Content.Find.Execute("61×14=", $true, $false, $false, $false, $false, $true, 1, $false, "17×62=", 2) | Out-Null
$d.Content.Find.Execute("26×89=", $true, $false, $false, $false, $false, $true, 1, $false, "80×96=", 2) | Out-Null
$d.Content.Find.Execute("82×35=", $true, $false, $false, $false, $false, $true, 1, $false, "68×38=", 2) | Out-Null
$d.Content.Find.Execute("35×53=", $true, $false, $false, $false, $false, $true, 1, $false, "38×82=", 2) | Out-Null
$d.Content.Find.Execute("22×41=", $true, $false, $false, $false, $false, $true, 1, $false, "93×89=", 2) | Out-Null
$d.Content.Find.Execute("13×19=", $true, $false, $false, $false, $false, $true, 1, $false, "75×37=", 2) | Out-Null
$d.Content.Find.Execute("57×92=", $true, $false, $false, $false, $false, $true, 1, $false, "37×21=", 2) | Out-Null
$d.Content.Find.Execute("42×87=", $true, $false, $false, $false, $false, $true, 1, $false, "38×43=", 2) | Out-Null
$d.Content.Find.Execute("71×87=", $true, $false, $false, $false, $false, $true, 1, $false, "91×35=", 2) | Out-Null
$d.Content.Find.Execute("76×13=", $true, $false, $false, $false, $false, $true, 1, $false, "69×12=", 2) | Out-Null
$d.Content.Find.Execute("72×90=", $true, $false, $false, $false, $false, $true, 1, $false, "94×25=", 2) | Out-Null
$d.Content.Find.Execute("90×15=", $true, $false, $false, $false, $false, $true, 1, $false, "18×20=", 2) | Out-Null
$d.Content.Find.Execute("39×100=", $true, $false, $false, $false, $false, $true, 1, $false, "42×52=", 2) | Out-Null
$d.Content.Find.Execute("68×91=", $true, $false, $false, $false, $false, $true, 1, $false, "67×66=", 2) | Out-Null
$d.Content.Find.Execute("70×97=", $true, $false, $false, $false, $false, $true, 1, $false, "53×72=", 2) | Out-Null
$d.Content.Find.Execute("77×67=", $true, $false, $false, $false, $false, $true, 1, $false, "15×46=", 2) | Out-Null
$d.Content.Find.Execute("71×73=", $true, $false, $false, $false, $false, $true, 1, $false, "34×25=", 2) | Out-Null
$d.Content.Find.Execute("60×18=", $true, $false, $false, $false, $false, $true, 1, $false, "33×31=", 2) | Out-Null
$d.Content.Find.Execute("65×52=", $true, $false, $false, $false, $false, $true, 1, $false, "78×42=", 2) | Out-Null
$d.Content.Find.Execute("42×49=", $true, $false, $false, $false, $false, $true, 1, $false, "96×95=", 2) | Out-Null
$d.Content.Find.Execute("48×69=", $true, $false, $false, $false, $false, $true, 1, $false, "69×55=", 2) | Out-Null
$d.Content.Find.Execute("41×21=", $true, $false, $false, $false, $false, $true, 1, $false, "35×23=", 2) | Out-Null
$d.Content.Find.Execute("67×95=", $true, $false, $false, $false, $false, $true, 1, $false, "82×56=", 2) | Out-Null
$d.Content.Find.Execute("27×24=", $true, $false, $false, $false, $false, $true, 1, $false, "16×47=", 2) | Out-Null
$d.Content.Find.Execute("87×85=", $true, $false, $false, $false, $false, $true, 1, $false, "64×43=", 2) | Out-Null
$d.Content.Find.Execute("36×10=", $true, $false, $false, $false, $false, $true, 1, $false, "98×23=", 2) | Out-Null
$d.Content.Find.Execute("56×15=", $true, $false, $false, $false, $false, $true, 1, $false, "59×96=", 2) | Out-Null
$d.Content.Find.Execute("86×97=", $true, $false, $false, $false, $false, $true, 1, $false, "60×35=", 2) | Out-Null
$d.Content.Find.Execute("74×15=", $true, $false, $false, $false, $false, $true, 1, $false, "52×80=", 2) | Out-Null
$d.Content.Find.Execute("47×35=", $true, $false, $false, $false, $false, $true, 1, $false, "66×23=", 2) | Out-Null
$d.Content.Find.Execute("63×79=", $true, $false, $false, $false, $false, $true, 1, $false, "23×67=", 2) | Out-Null
$d.Content.Find.Execute("57×49=", $true, $false, $false, $false, $false, $true, 1, $false, "18×48=", 2) | Out-Null
$d.Content.Find.Execute("22×83=", $true, $false, $false, $false, $false, $true, 1, $false, "86×46=", 2) | Out-Null
$d.Content.Find.Execute("29×72=", $true, $false, $false, $false, $false, $true, 1, $false, "15×38=", 2) | Out-Null
$d.Content.Find.Execute("35×66=", $true, $false, $false, $false, $false, $true, 1, $false, "25×51=", 2) | Out-Null
$d.Content.Find.Execute("17×28=", $true, $false, $false, $false, $false, $true, 1, $false, "44×76=", 2) | Out-Null
$d.Content.Find.Execute("43×52=", $true, $false, $false, $false, $false, $true, 1, $false, "40×47=", 2) | Out-Null
$d.Content.Find.Execute("23×21=", $true, $false, $false, $false, $false, $true, 1, $false, "28×62=", 2) | Out-Null
$d.Content.Find.Execute("22×45=", $true, $false, $false, $false, $false, $true, 1, $false, "61×45=", 2) | Out-Null
$d.Content.Find.Execute("25×24=", $true, $false, $false, $false, $false, $true, 1, $false, "70×33=", 2) | Out-Null
$d.Content.Find.Execute("73×16=", $true, $false, $false, $false, $false, $true, 1, $false, "95×96=", 2) | Out-Null
$d.Content.Find.Execute("24×74=", $true, $false, $false, $false, $false, $true, 1, $false, "85×59=", 2) | Out-Null
